$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162 (pushes existing rows 162-171 down to 163-172,
# carrying their data/formatting with them).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly price-report entry.
$ws.Cells.Item(162, 1).Value = 8
$ws.Cells.Item(162, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(162, 3).Value = "Coquimbo"
$ws.Cells.Item(162, 4).Value = 44516
$ws.Cells.Item(162, 5).Value = 4
$ws.Cells.Item(162, 6).Value = 100112003
$ws.Cells.Item(162, 7).Value = "Ajo"
$ws.Cells.Item(162, 8).Value = "Chino"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 500
$ws.Cells.Item(162, 11).Value = 18500
$ws.Cells.Item(162, 12).Value = 19000
$ws.Cells.Item(162, 13).Value = 18750
$ws.Cells.Item(162, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(162, 15).Value = "China"
$ws.Cells.Item(162, 16).Value = 1875
$ws.Cells.Item(162, 17).Value = 10
$ws.Cells.Item(162, 18).Value = "Hortaliza"
